$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) for rows 2-27, columns B-E and G (TB, d2S, K, IP, sum)
$ws.Cells.Item(2, 2).Value = 0.7287194209349384; $ws.Cells.Item(2, 3).Value = 0.3375848360084654; $ws.Cells.Item(2, 4).Value = 0.1529057820181812; $ws.Cells.Item(2, 5).Value = 0.4998867070740569; $ws.Cells.Item(2, 7).Value = 1.719096746035642
$ws.Cells.Item(3, 2).Value = 3.182878228561681; $ws.Cells.Item(3, 3).Value = 1.65323645889881; $ws.Cells.Item(3, 4).Value = 16.98373111632243; $ws.Cells.Item(3, 5).Value = 0.4998867070740569; $ws.Cells.Item(3, 7).Value = 22.31973251085698
$ws.Cells.Item(4, 2).Value = 0.1554434735375247; $ws.Cells.Item(4, 3).Value = 1.65323645889881; $ws.Cells.Item(4, 4).Value = 0.1529057820181812; $ws.Cells.Item(4, 5).Value = 0.4998867070740569; $ws.Cells.Item(4, 7).Value = 2.461472421528573
$ws.Cells.Item(5, 2).Value = 1.505614041169197; $ws.Cells.Item(5, 3).Value = 1.65323645889881; $ws.Cells.Item(5, 4).Value = 0.1529057820181812; $ws.Cells.Item(5, 5).Value = 0.4998867070740569; $ws.Cells.Item(5, 7).Value = 3.811642989160245
$ws.Cells.Item(6, 2).Value = 3.182878228561681; $ws.Cells.Item(6, 3).Value = 1.65323645889881; $ws.Cells.Item(6, 4).Value = 0.7127328510149897; $ws.Cells.Item(6, 5).Value = 0.4998867070740569; $ws.Cells.Item(6, 7).Value = 6.048734245549538
$ws.Cells.Item(7, 2).Value = 3.182878228561681; $ws.Cells.Item(7, 3).Value = 1.65323645889881; $ws.Cells.Item(7, 4).Value = 0.1529057820181812; $ws.Cells.Item(7, 5).Value = 0.4998867070740569; $ws.Cells.Item(7, 7).Value = 5.488907176552729
$ws.Cells.Item(8, 2).Value = 0.02258322285507441; $ws.Cells.Item(8, 3).Value = 0.004309184025731883; $ws.Cells.Item(8, 4).Value = 0.1529057820181812; $ws.Cells.Item(8, 5).Value = 0.4998867070740569; $ws.Cells.Item(8, 7).Value = 0.6796848959730444
$ws.Cells.Item(9, 2).Value = 3.182878228561681; $ws.Cells.Item(9, 3).Value = 1.65323645889881; $ws.Cells.Item(9, 4).Value = 0.1529057820181812; $ws.Cells.Item(9, 5).Value = 0.4998867070740569; $ws.Cells.Item(9, 7).Value = 5.488907176552729
$ws.Cells.Item(10, 2).Value = 0.1554434735375247; $ws.Cells.Item(10, 3).Value = 0.3375848360084654; $ws.Cells.Item(10, 4).Value = 0.1529057820181812; $ws.Cells.Item(10, 5).Value = 0.4998867070740569; $ws.Cells.Item(10, 7).Value = 1.145820798638228
$ws.Cells.Item(11, 2).Value = 1.505614041169197; $ws.Cells.Item(11, 3).Value = 0.3375848360084654; $ws.Cells.Item(11, 4).Value = 0.7127328510149897; $ws.Cells.Item(11, 5).Value = 0.4998867070740569; $ws.Cells.Item(11, 7).Value = 3.055818435266709
$ws.Cells.Item(12, 2).Value = 3.182878228561681; $ws.Cells.Item(12, 3).Value = 1.65323645889881; $ws.Cells.Item(12, 4).Value = 0.1529057820181812; $ws.Cells.Item(12, 5).Value = 0.4998867070740569; $ws.Cells.Item(12, 7).Value = 5.488907176552729
$ws.Cells.Item(13, 2).Value = 1.505614041169197; $ws.Cells.Item(13, 3).Value = 1.65323645889881; $ws.Cells.Item(13, 4).Value = 3.082599426703578; $ws.Cells.Item(13, 5).Value = 0.4998867070740569; $ws.Cells.Item(13, 7).Value = 6.741336633845642
$ws.Cells.Item(14, 2).Value = 0.7287194209349384; $ws.Cells.Item(14, 3).Value = 1.65323645889881; $ws.Cells.Item(14, 4).Value = 3.082599426703578; $ws.Cells.Item(14, 5).Value = 0.4998867070740569; $ws.Cells.Item(14, 7).Value = 5.964442013611383
$ws.Cells.Item(15, 2).Value = 3.182878228561681; $ws.Cells.Item(15, 3).Value = 1.65323645889881; $ws.Cells.Item(15, 4).Value = 16.98373111632243; $ws.Cells.Item(15, 5).Value = 0.4998867070740569; $ws.Cells.Item(15, 7).Value = 22.31973251085698
$ws.Cells.Item(16, 2).Value = 0.06328177979961902; $ws.Cells.Item(16, 3).Value = 0.3375848360084654; $ws.Cells.Item(16, 4).Value = 0.7127328510149897; $ws.Cells.Item(16, 5).Value = 0.4998867070740569; $ws.Cells.Item(16, 7).Value = 1.613486173897131
$ws.Cells.Item(17, 2).Value = 3.182878228561681; $ws.Cells.Item(17, 3).Value = 1.65323645889881; $ws.Cells.Item(17, 4).Value = 3.082599426703578; $ws.Cells.Item(17, 5).Value = 0.4998867070740569; $ws.Cells.Item(17, 7).Value = 8.418600821238126
$ws.Cells.Item(18, 2).Value = 0.7287194209349384; $ws.Cells.Item(18, 3).Value = 0.3375848360084654; $ws.Cells.Item(18, 4).Value = 0.7127328510149897; $ws.Cells.Item(18, 5).Value = 0.4998867070740569; $ws.Cells.Item(18, 7).Value = 2.27892381503245
$ws.Cells.Item(19, 2).Value = 3.182878228561681; $ws.Cells.Item(19, 3).Value = 1.65323645889881; $ws.Cells.Item(19, 4).Value = 0.7127328510149897; $ws.Cells.Item(19, 5).Value = 0.4998867070740569; $ws.Cells.Item(19, 7).Value = 6.048734245549538
$ws.Cells.Item(20, 2).Value = 3.182878228561681; $ws.Cells.Item(20, 3).Value = 1.65323645889881; $ws.Cells.Item(20, 4).Value = 0.7127328510149897; $ws.Cells.Item(20, 5).Value = 0.4998867070740569; $ws.Cells.Item(20, 7).Value = 6.048734245549538
$ws.Cells.Item(21, 2).Value = 3.182878228561681; $ws.Cells.Item(21, 3).Value = 1.65323645889881; $ws.Cells.Item(21, 4).Value = 16.98373111632243; $ws.Cells.Item(21, 5).Value = 6.48142807727062; $ws.Cells.Item(21, 7).Value = 28.30127388105354
$ws.Cells.Item(22, 2).Value = 1.505614041169197; $ws.Cells.Item(22, 3).Value = 1.65323645889881; $ws.Cells.Item(22, 4).Value = 0.7127328510149897; $ws.Cells.Item(22, 5).Value = 0.4998867070740569; $ws.Cells.Item(22, 7).Value = 4.371470058157054
$ws.Cells.Item(23, 2).Value = 3.182878228561681; $ws.Cells.Item(23, 3).Value = 9.226618575922256; $ws.Cells.Item(23, 4).Value = 3.082599426703578; $ws.Cells.Item(23, 5).Value = 6.48142807727062; $ws.Cells.Item(23, 7).Value = 21.97352430845813
$ws.Cells.Item(24, 2).Value = 3.182878228561681; $ws.Cells.Item(24, 3).Value = 1.65323645889881; $ws.Cells.Item(24, 4).Value = 3.082599426703578; $ws.Cells.Item(24, 5).Value = 0.4998867070740569; $ws.Cells.Item(24, 7).Value = 8.418600821238126
$ws.Cells.Item(25, 2).Value = 3.182878228561681; $ws.Cells.Item(25, 3).Value = 1.65323645889881; $ws.Cells.Item(25, 4).Value = 3.082599426703578; $ws.Cells.Item(25, 5).Value = 0.4998867070740569; $ws.Cells.Item(25, 7).Value = 8.418600821238126
$ws.Cells.Item(26, 2).Value = 3.182878228561681; $ws.Cells.Item(26, 3).Value = 1.65323645889881; $ws.Cells.Item(26, 4).Value = 0.1529057820181812; $ws.Cells.Item(26, 5).Value = 0.4998867070740569; $ws.Cells.Item(26, 7).Value = 5.488907176552729
$ws.Cells.Item(27, 2).Value = 3.182878228561681; $ws.Cells.Item(27, 3).Value = 1.65323645889881; $ws.Cells.Item(27, 4).Value = 0.7127328510149897; $ws.Cells.Item(27, 5).Value = 0.4998867070740569; $ws.Cells.Item(27, 7).Value = 6.048734245549538
